$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 6-10, and add a new row 11, with text values
# matching the shared-string (t="s") storage of the original cells.
# Force text number format so numeric-looking strings stay text (t="s"),
# matching the original workbook's storage convention.

$textRange = $ws.Range("A6:C11")
$textRange.NumberFormat = "@"

$ws.Range("C6").Value = "9"
$ws.Range("C7").Value = "8"

$ws.Range("A9").Value = "5"
$ws.Range("B9").Value = "Stand de los Besos 3"
$ws.Range("C9").Value = "3"

$ws.Range("A10").Value = "6"
$ws.Range("B10").Value = "La Casa de Papel Temporada 1"
$ws.Range("C10").Value = "2"

$ws.Range("A11").Value = "4"
$ws.Range("B11").Value = "Stand de los Besos 2"
$ws.Range("C11").Value = "2"
